$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "52.119.11"
$ws.Range("D2").ClearFormats()
$ws.Range("E2").Value = "  +0.85%  "

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.941.54"
$ws.Range("D3").ClearFormats()
$ws.Range("E3").Value = "  +4.58%  "

# Row 4
$ws.Range("E4").Value = "  +0.02%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "353.40"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  +0.75%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "112.17"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  -0.82%  "

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.563"
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = "  +0.49%  "

# Row 8
$ws.Range("E8").Value = "  +0.00%  "

# Row 9
$ws.Range("E9").Value = "  +0.89%  "

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "39.51"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "  -1.41%  "

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0879"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = "  +3.51%  "

# Row 12
$ws.Range("E12").Value = "  +1.08%  "

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "20.15"
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = "  +0.85%  "

# Row 14
$ws.Range("B14").Value = "Polkadot"
$ws.Range("C14").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "7.77"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "  +0.10%  "

# Row 15
$ws.Range("B15").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C15").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "3.401.80"
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = "  +4.43%  "

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "2.947.20"
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = "  +4.80%  "

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.983"
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = "  +1.03%  "

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "52.163.43"
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = "  +0.81%  "

# Row 19
$ws.Range("E19").Value = "  +0.50%  "

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "3.30"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "  -3.31%  "

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "14.25"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "  +6.16%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.0₃0981"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  +0.99%  "

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "71.26"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "  +1.07%  "

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "268.55"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "  -0.04%  "

# Row 25
$ws.Range("E25").Value = "  +1.14%  "

# Row 26
$ws.Range("E26").Value = "  +10.88%  "

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "27.06"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "  +3.28%  "

# Row 28
$ws.Range("E28").Value = "  -0.06%  "

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "7.31"
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = "  +18.38%  "

# Row 30
$ws.Range("E30").Value = "  +16.09%  "

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "10.62"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "  +0.72%  "

# Row 32
$ws.Range("E32").Value = "  -0.31%  "

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "37.12"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = "  -4.16%  "

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "6.11"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = "  +7.74%  "

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "53.06"
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = "  +0.66%  "

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.0453"
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = "  +0.02%  "

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.998"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "  -0.13%  "

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "3.39"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = "  +6.01%  "

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "18.70"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "  -2.41%  "

# Row 40
$ws.Range("E40").Value = "  +2.68%  "

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "2.70"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "  +5.28%  "

# Row 42
$ws.Range("E42").Value = "  +1.71%  "

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "23.31"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  +5.34%  "

# Row 44
$ws.Range("E44").Value = "  -1.53%  "

# Row 45
$ws.Range("B45").Value = "Maker"
$ws.Range("C45").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "2.201.79"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "  +2.22%  "

# Row 46
$ws.Range("B46").Value = "ApeXProtocol"
$ws.Range("C46").Value = "https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.53"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "  +1.72%  "

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "3.53"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "  +0.96%  "

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "111.25"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "  -8.60%  "

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.247"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = "  +9.92%  "

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0355"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "  +10.62%  "

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.956"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  -2.97%  "
